# Updated cryptos list on Sat May 18 17:36:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: column letter -> new text value.
# (Row 46/47 also swap Coin name + Link, since VeChain/dogwifhat traded
# ranking positions in this refresh.)
$updates = @{
    2  = @{ D = "66.870.19";  E = "  -0.03%  " }
    3  = @{ D = "3.115.84";   E = "  +0.71%  " }
    4  = @{                   E = "  +0.01%  " }
    5  = @{ D = "577.64";     E = "  -0.43%  " }
    6  = @{ D = "171.15";     E = "  +1.98%  " }
    7  = @{                   E = "  -0.03%  " }
    8  = @{ D = "3.113.85";   E = "  +0.80%  " }
    9  = @{ D = "0.521";      E = "  -0.74%  " }
    10 = @{ D = "6.47";       E = "  -2.80%  " }
    11 = @{                   E = "  -1.37%  " }
    12 = @{                   E = "  +0.23%  " }
    13 = @{ D = "0.0000245";  E = "  -1.97%  " }
    14 = @{ D = "37.20";      E = "  +1.31%  " }
    15 = @{                   E = "  -1.13%  " }
    16 = @{ D = "3.634.66";   E = "  +0.65%  " }
    17 = @{ D = "66.970.18";  E = "  +0.05%  " }
    18 = @{ D = "7.17";       E = "  -0.87%  " }
    19 = @{ D = "3.115.41";   E = "  +0.49%  " }
    20 = @{ D = "16.31";      E = "  -0.44%  " }
    21 = @{ D = "476.77";     E = "  +1.55%  " }
    22 = @{                   E = "  -0.17%  " }
    23 = @{ D = "7.93";       E = "  +4.92%  " }
    24 = @{ D = "13.42";      E = "  +4.37%  " }
    25 = @{ D = "84.06";      E = "  +0.93%  " }
    26 = @{                   E = "  -3.12%  " }
    27 = @{ D = "10.11";      E = "  -0.51%  " }
    28 = @{                   E = "  +0.03%  " }
    29 = @{                   E = "  -2.22%  " }
    30 = @{ D = "2.38";       E = "  -1.76%  " }
    31 = @{ D = "2.67";       E = "  -0.21%  " }
    32 = @{ D = "28.52";      E = "  +1.05%  " }
    33 = @{                   E = "  +0.18%  " }
    34 = @{ D = "0.0₃0939";   E = "  -8.73%  " }
    35 = @{                   E = "  +0.00%  " }
    36 = @{                   E = "  -0.65%  " }
    37 = @{ D = "0.972";      E = "  -3.35%  " }
    38 = @{ D = "46.91";      E = "  +0.29%  " }
    39 = @{ D = "2.07";       E = "  -0.51%  " }
    40 = @{                   E = "  -0.46%  " }
    41 = @{                   E = "  -1.86%  " }
    42 = @{                   E = "  -0.93%  " }
    43 = @{ D = "8.72";       E = "  +0.01%  " }
    44 = @{ D = "2.854.72";   E = "  +3.40%  " }
    45 = @{ D = "387.69";     E = "  -0.07%  " }
    46 = @{ B = "VeChain";   C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.0358"; E = "  -1.54%  " }
    47 = @{ B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif";    D = "2.59";   E = "  -9.05%  " }
    48 = @{ D = "135.97";     E = "  +0.82%  " }
    50 = @{ D = "24.89";      E = "  +0.32%  " }
    51 = @{                   E = "  -2.02%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $newVal = $cols[$col]

        if ($col -eq "D") {
            # D holds plain price text (possibly using "." as a thousands
            # separator, e.g. "66.870.19") that must stay text. A handful of
            # the new values (e.g. "577.64") look like plain decimals and
            # Excel's Range.Value setter would silently coerce them to a
            # Number. Force text entry by briefly marking the cell as Text,
            # then restore the plain "Normal" cell style so no formatting
            # residue is left behind on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newVal
        }
    }
}
